$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.214.42"
$ws.Range("E2").Value = "  -3.70%  "
$ws.Range("D3").Value = "'3.506.35"
$ws.Range("E3").Value = "  -5.56%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'582.38"
$ws.Range("E5").Value = "  -1.24%  "
$ws.Range("D6").Value = "'173.33"
$ws.Range("E6").Value = "  -3.90%  "
$ws.Range("D7").Value = "'0.621"
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("D8").Value = "'3.500.77"
$ws.Range("E8").Value = "  -5.51%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "'0.189"
$ws.Range("E10").Value = "  -6.96%  "
$ws.Range("D11").Value = "'6.74"
$ws.Range("E11").Value = "  +4.34%  "
$ws.Range("D12").Value = "'0.595"
$ws.Range("E12").Value = "  -3.11%  "
$ws.Range("D13").Value = "'46.85"
$ws.Range("E13").Value = "  -6.40%  "
$ws.Range("D14").Value = "'0.0000276"
$ws.Range("E14").Value = "  -4.19%  "
$ws.Range("D15").Value = "'676.94"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("D16").Value = "'4.078.66"
$ws.Range("E16").Value = "  -5.43%  "
$ws.Range("D17").Value = "'8.71"
$ws.Range("E17").Value = "  -3.74%  "
$ws.Range("D18").Value = "'69.272.59"
$ws.Range("E18").Value = "  -3.72%  "
$ws.Range("D19").Value = "'3.518.47"
$ws.Range("E19").Value = "  -5.32%  "
$ws.Range("E20").Value = "  -1.35%  "
$ws.Range("D21").Value = "'17.44"
$ws.Range("E21").Value = "  -3.90%  "
$ws.Range("D22").Value = "'11.18"
$ws.Range("E22").Value = "  -4.42%  "
$ws.Range("D23").Value = "'0.901"
$ws.Range("E23").Value = "  -4.70%  "
$ws.Range("D24").Value = "'16.16"
$ws.Range("E24").Value = "  -9.49%  "
$ws.Range("D25").Value = "'97.96"
$ws.Range("E25").Value = "  -5.71%  "
$ws.Range("D26").Value = "'3.87"
$ws.Range("E26").Value = "  -4.73%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "'2.66"
$ws.Range("E28").Value = "  -6.53%  "
$ws.Range("D29").Value = "'9.40"
$ws.Range("E29").Value = "  -8.99%  "
$ws.Range("D30").Value = "'32.97"
$ws.Range("E30").Value = "  -7.37%  "
$ws.Range("D31").Value = "'8.72"
$ws.Range("E31").Value = "  -6.21%  "
$ws.Range("D32").Value = "'3.19"
$ws.Range("E32").Value = "  -7.99%  "
$ws.Range("B33").Value = "Mantle"
$ws.Range("C33").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D33").Value = "'1.36"
$ws.Range("E33").Value = "  -6.13%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "'7.28"
$ws.Range("E34").Value = "  -1.29%  "
$ws.Range("D35").Value = "'597.64"
$ws.Range("E35").Value = "  +6.34%  "
$ws.Range("D36").Value = "'3.59"
$ws.Range("E36").Value = "  -15.89%  "
$ws.Range("D37").Value = "'10.87"
$ws.Range("E37").Value = "  -3.86%  "
$ws.Range("E38").Value = "  -5.08%  "
$ws.Range("D39").Value = "'57.29"
$ws.Range("E39").Value = "  -4.00%  "
$ws.Range("E40").Value = "  +0.22%  "
$ws.Range("D41").Value = "'0.0438"
$ws.Range("E41").Value = "  -5.66%  "
$ws.Range("D42").Value = "'0.336"
$ws.Range("E42").Value = "  -5.05%  "
$ws.Range("D43").Value = "'3.418.53"
$ws.Range("E43").Value = "  -9.45%  "
$ws.Range("E44").Value = "  -6.14%  "
$ws.Range("D45").Value = "'33.33"
$ws.Range("E45").Value = "  -6.93%  "
$ws.Range("D46").Value = "0.0₃0707"
$ws.Range("E46").Value = "  -9.36%  "
$ws.Range("D47").Value = "'2.90"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("D48").Value = "'2.60"
$ws.Range("E48").Value = "  -7.52%  "
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("D50").Value = "'5.78"
$ws.Range("E50").Value = "  +17.53%  "
$ws.Range("E51").Value = "  -2.15%  "
